$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.461.21"
$ws.Range("E2").Value = "  +1.09%  "
$ws.Range("D3").Value = "2.332.42"
$ws.Range("E3").Value = "  +1.24%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "303.30"
$ws.Range("E5").Value = "  +1.10%  "
$ws.Range("D6").Value = "98.42"
$ws.Range("E6").Value = "  +0.68%  "
$ws.Range("E7").Value = "  -0.69%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -0.18%  "
$ws.Range("D10").Value = "35.83"
$ws.Range("E10").Value = "  -0.30%  "
$ws.Range("D11").Value = "19.57"
$ws.Range("E11").Value = "  +7.51%  "
$ws.Range("D12").Value = "0.0801"
$ws.Range("E12").Value = "  +1.31%  "
$ws.Range("E13").Value = "  +0.28%  "
$ws.Range("E14").Value = "  +2.05%  "
$ws.Range("D15").Value = "2.693.42"
$ws.Range("E15").Value = "  +1.09%  "
$ws.Range("D16").Value = "2.317.66"
$ws.Range("E16").Value = "  +0.72%  "
$ws.Range("D17").Value = "0.793"
$ws.Range("E17").Value = "  +1.49%  "
$ws.Range("D18").Value = "43.377.93"
$ws.Range("E18").Value = "  +1.07%  "
$ws.Range("D19").Value = "12.79"
$ws.Range("E19").Value = "  +1.06%  "
$ws.Range("D20").Value = "0.0₃0903"
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("E21").Value = "  +1.12%  "
$ws.Range("D22").Value = "68.02"
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("E23").Value = "  +0.89%  "
$ws.Range("D24").Value = "2.24"
$ws.Range("E24").Value = "  +4.51%  "
$ws.Range("D25").Value = "2.45"
$ws.Range("E25").Value = "  +0.24%  "
$ws.Range("E26").Value = "  +0.20%  "
$ws.Range("D27").Value = "25.13"
$ws.Range("E27").Value = "  -1.22%  "
$ws.Range("D28").Value = "2.21"
$ws.Range("E28").Value = "  +7.71%  "
$ws.Range("D29").Value = "164.97"
$ws.Range("E29").Value = "  -0.20%  "
$ws.Range("E30").Value = "  +1.03%  "
$ws.Range("D31").Value = "33.50"
$ws.Range("E31").Value = "  +0.97%  "
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("E33").Value = "  +0.24%  "
$ws.Range("D34").Value = "17.97"
$ws.Range("E34").Value = "  +6.03%  "
$ws.Range("D35").Value = "4.51"
$ws.Range("E35").Value = "  -6.84%  "
$ws.Range("E36").Value = "  +2.47%  "
$ws.Range("E37").Value = "  -1.38%  "
$ws.Range("E38").Value = "  -0.20%  "
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").Value = "1.78"
$ws.Range("E39").Value = "  +0.86%  "
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").Value = "2.80"
$ws.Range("E40").Value = "  +2.19%  "
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("D42").Value = "1.992.73"
$ws.Range("E42").Value = "  -0.81%  "
$ws.Range("D43").Value = "10.79"
$ws.Range("E43").Value = "  +7.18%  "
$ws.Range("D44").Value = "0.0282"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").Value = "18.18"
$ws.Range("E45").Value = "  +3.38%  "
$ws.Range("D46").Value = "2.06"
$ws.Range("E46").Value = "  -1.48%  "
$ws.Range("D47").Value = "2.81"
$ws.Range("E47").Value = "  +0.70%  "
$ws.Range("B48").Value = "MultiversX"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D48").Value = "54.24"
$ws.Range("E48").Value = "  +1.23%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "2.561.52"
$ws.Range("E49").Value = "  +1.18%  "
$ws.Range("B50").Value = "HuobiToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D50").Value = "2.86"
$ws.Range("E50").Value = "  -3.59%  "
$ws.Range("D51").Value = "72.88"
$ws.Range("E51").Value = "  +0.87%  "
